$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, [string]$value)
    # Force the cell to remain a text value (matches the source file's
    # inline-string cells) instead of Excel auto-coercing numeric-looking
    # strings ("529.20", "1.00", ...) into Number cells and dropping
    # trailing zeros / changing the stored type. Resetting the style back
    # to Normal afterwards avoids leaving a stray text number-format on
    # the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell $ws.Range("D2") "59.122.90"
$ws.Range("E2").Value = "  +1.48%  "

# Row 3 - Ethereum
Set-TextCell $ws.Range("D3") "2.591.20"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.06%  "

# Row 5 - BNB
Set-TextCell $ws.Range("D5") "529.20"
$ws.Range("E5").Value = "  +2.03%  "

# Row 6 - Solana
Set-TextCell $ws.Range("D6") "139.92"
$ws.Range("E6").Value = "  -0.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - XRP
Set-TextCell $ws.Range("D8") "0.566"
$ws.Range("E8").Value = "  +0.51%  "

# Row 9 - LidoStakedEther
Set-TextCell $ws.Range("D9") "2.603.36"
$ws.Range("E9").Value = "  +0.00%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -0.39%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.40%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +1.09%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +3.08%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell $ws.Range("D14") "3.055.03"
$ws.Range("E14").Value = "  +0.43%  "

# Row 15 - WrappedBTC
Set-TextCell $ws.Range("D15") "59.053.77"
$ws.Range("E15").Value = "  +1.41%  "

# Row 16 - Avalanche
$ws.Range("E16").Value = "  +1.01%  "

# Row 17 - ShibaInu
Set-TextCell $ws.Range("D17") "0.0000134"
$ws.Range("E17").Value = "  +0.48%  "

# Row 18 - WrappedEther
Set-TextCell $ws.Range("D18") "2.589.16"
$ws.Range("E18").Value = "  +0.14%  "

# Row 19 - BitcoinCash
Set-TextCell $ws.Range("D19") "347.60"
$ws.Range("E19").Value = "  +3.23%  "

# Row 20 - Polkadot
Set-TextCell $ws.Range("D20") "4.34"
$ws.Range("E20").Value = "  +0.80%  "

# Row 21 - Chainlink
Set-TextCell $ws.Range("D21") "10.09"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22 - Uniswap
Set-TextCell $ws.Range("D22") "6.40"
$ws.Range("E22").Value = "  +0.22%  "

# Row 23 - Dai
Set-TextCell $ws.Range("D23") "0.998"
$ws.Range("E23").Value = "  +0.13%  "

# Row 24 - Litecoin
Set-TextCell $ws.Range("D24") "67.48"
$ws.Range("E24").Value = "  +2.70%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -0.39%  "

# Row 26 - Polygon
$ws.Range("E26").Value = "  +1.20%  "

# Row 27 - Binance-PegBSC-USD
Set-TextCell $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.30%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("E28").Value = "  +2.33%  "

# Row 29 - USDe
$ws.Range("E29").Value = "  +0.06%  "

# Row 30 - PEPE
Set-TextCell $ws.Range("D30") "0.0₃0734"
$ws.Range("E30").Value = "  +0.16%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.74%  "

# Row 32 - Aptos
Set-TextCell $ws.Range("D32") "5.86"
$ws.Range("E32").Value = "  -3.14%  "

# Row 33 - EthereumClassic
Set-TextCell $ws.Range("D33") "18.78"
$ws.Range("E33").Value = "  +0.35%  "

# Row 34 - Monero
Set-TextCell $ws.Range("D34") "148.72"
$ws.Range("E34").Value = "  -0.27%  "

# Row 35 - NEARProtocol
$ws.Range("E35").Value = "  +0.60%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  -0.41%  "

# Row 37 - OKB
Set-TextCell $ws.Range("D37") "36.90"
$ws.Range("E37").Value = "  +1.97%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +0.61%  "

# Row 39 - SuiNetwork
Set-TextCell $ws.Range("D39") "0.827"
$ws.Range("E39").Value = "  +0.21%  "

# Row 40 - Fetch.AI
Set-TextCell $ws.Range("D40") "0.829"
$ws.Range("E40").Value = "  -2.61%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +0.83%  "

# Row 42 - FirstDigitalUSD
Set-TextCell $ws.Range("D42") "0.997"
$ws.Range("E42").Value = "  +0.10%  "

# Row 43 - now Bittensor (was WhiteBITCoin)
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell $ws.Range("D43") "269.23"
$ws.Range("E43").Value = "  -1.36%  "

# Row 44 - now WhiteBITCoin (was Bittensor)
$ws.Range("B44").Value = "WhiteBITCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell $ws.Range("D44") "10.74"
$ws.Range("E44").Value = "  +0.50%  "

# Row 45 - Mantle
Set-TextCell $ws.Range("D45") "0.595"
$ws.Range("E45").Value = "  -1.10%  "

# Row 46 - Stellar
$ws.Range("E46").Value = "  +1.23%  "

# Row 47 - Hedera
$ws.Range("E47").Value = "  -0.10%  "

# Row 48 - EnergySwap
Set-TextCell $ws.Range("D48") "18.41"
$ws.Range("E48").Value = "  -0.67%  "

# Row 49 - Maker
Set-TextCell $ws.Range("D49") "1.951.33"
$ws.Range("E49").Value = "  -0.86%  "

# Row 50 - now VeChain (was RenderToken)
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws.Range("D50") "0.0221"
$ws.Range("E50").Value = "  +0.44%  "

# Row 51 - now InjectiveProtocol (was VeChain)
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws.Range("D51") "18.18"
$ws.Range("E51").Value = "  +0.28%  "
